$wb = $excel.ActiveWorkbook

# --- Insert the new "Reduced" sheet between "ST muscles" and "reduced muscles" ---
$stMuscles = $wb.Worksheets.Item("ST muscles")
$reducedMuscles = $wb.Worksheets.Item("reduced muscles")
$newSheet = $wb.Worksheets.Add($reducedMuscles)
$newSheet.Name = "Reduced"

Write-Output "sheet added"
